$wb = $excel.ActiveWorkbook

# --- Add the new "Tips" worksheet as the last sheet in the workbook ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$ws.Name = 'Tips'

# --- Column widths (approximate Excel's auto-fit pixel widths) ---
$ws.Columns("B").ColumnWidth = 31.6
$ws.Columns("C").ColumnWidth = 89.5

# --- Row 2 ---
$ws.Range("B2").Value = 'backward \forward attribute selection'
$ws.Range("C2").Value = 'start from 1 column to train the system and then add more or vise versa'

# --- Row 3 ---
$ws.Range("B3").Value = 'Correlation analysis'
$ws.Range("C3").Value = "if two or more columns are corelated, for example for big houses land and size of the houses are always bigger than for small houses we don't need to keep land size and house size columns"

# --- Row 4 ---
$ws.Range("B4").Value = 'Datasets'
$ws.Range("C4").Value = 'https://archive.ics.uci.edu/ml/datasets.php'

# --- Row 5 ---
$ws.Range("C5").Value = 'https://www.kaggle.com/datasets'

# --- Row 7 ---
$ws.Range("B7").Value = 'Balanced Datframe'
$ws.Range("C7").Value = 'when target column have similar numbers of 1 and 0'

# --- Row 8 ---
$ws.Range("B8").Value = 'EDA - explarotatry Data analysis'
$ws.Range("C8").Value = 'learn your data: shape(), dtypes, isna().sum(), describe(), df.["target"].value_counts(), corr()'

# --- Alignment: column C wrap text, column B vertical-top (order matters: it
# --- determines the order new cell-style records are appended in styles.xml) ---
foreach ($addr in @('C2', 'C3', 'C7', 'C8')) {
    $ws.Range($addr).WrapText = $true
}
foreach ($addr in @('B2', 'B3', 'B4', 'B7', 'B8')) {
    $ws.Range($addr).VerticalAlignment = -4160
}

# --- Hyperlinks (also gives C4/C5 the Hyperlink+wrap style) ---
$ws.Range("C4").WrapText = $true
$ws.Range("C5").WrapText = $true
$ws.Hyperlinks.Add($ws.Range("C4"), 'https://archive.ics.uci.edu/ml/datasets.php') | Out-Null
$ws.Hyperlinks.Add($ws.Range("C5"), 'https://www.kaggle.com/datasets') | Out-Null

# --- Row height for the wrapped correlation-analysis row ---
$ws.Rows("3").RowHeight = 28.8

# --- Selection / active-tab bookkeeping ---
$ws.Range("C12").Select() | Out-Null
$ws.Activate()
